$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Ranking")
$ws.Range("H4").Value = 0.002542775759326063
$ws.Range("I4").Value = 0.002376293918721295
$ws.Range("H8").Value = 0.003817475084992743
$ws.Range("I8").Value = 0.002893976026713732

$ws = $wb.Worksheets.Item("Matriz_Pvalores")
$ws.Range("E2").Value = 0.003295502448879439
$ws.Range("J2").Value = 0.00315039942074824
$ws.Range("E3").Value = 0.123675267522446
$ws.Range("J3").Value = 0.07046843728652985
$ws.Range("E4").Value = 0.3889653978921017
$ws.Range("J4").Value = 0.05219548808690333
$ws.Range("B5").Value = 0.003295502448879439
$ws.Range("C5").Value = 0.123675267522446
$ws.Range("D5").Value = 0.3889653978921017
$ws.Range("F5").Value = 0.00288943787630247
$ws.Range("G5").Value = 0.30439985734797
$ws.Range("H5").Value = 0.2359497530506454
$ws.Range("I5").Value = 0.5249598689705572
$ws.Range("J5").Value = 0.2365482656164806
$ws.Range("E6").Value = 0.00288943787630247
$ws.Range("J6").Value = 0.0008798021143592294
$ws.Range("E7").Value = 0.30439985734797
$ws.Range("J7").Value = 0.8475117616320325
$ws.Range("E8").Value = 0.2359497530506454
$ws.Range("J8").Value = 0.6399236992372721
$ws.Range("E9").Value = 0.5249598689705572
$ws.Range("J9").Value = 0.5314153404758986
$ws.Range("B10").Value = 0.00315039942074824
$ws.Range("C10").Value = 0.07046843728652985
$ws.Range("D10").Value = 0.05219548808690333
$ws.Range("E10").Value = 0.2365482656164806
$ws.Range("F10").Value = 0.0008798021143592294
$ws.Range("G10").Value = 0.8475117616320325
$ws.Range("H10").Value = 0.6399236992372721
$ws.Range("I10").Value = 0.5314153404758986

$ws = $wb.Worksheets.Item("Matriz_DM_Original")
$ws.Range("E2").Value = 6.273444151869245
$ws.Range("J2").Value = 6.350371468051042
$ws.Range("E3").Value = -1.944928082375265
$ws.Range("J3").Value = -2.449694784569639
$ws.Range("E4").Value = 0.9655088164592457
$ws.Range("J4").Value = 2.734589956512321
$ws.Range("B5").Value = -6.273444151869245
$ws.Range("C5").Value = 1.944928082375265
$ws.Range("D5").Value = -0.9655088164592457
$ws.Range("F5").Value = -6.500345077603687
$ws.Range("G5").Value = 1.177142446357576
$ws.Range("H5").Value = 1.393344613414188
$ws.Range("I5").Value = 0.6956293248699703
$ws.Range("J5").Value = 1.391201525421976
$ws.Range("E6").Value = 6.500345077603687
$ws.Range("J6").Value = 8.902724452679363
$ws.Range("E7").Value = -1.177142446357576
$ws.Range("J7").Value = 0.2050953498393789
$ws.Range("E8").Value = -1.393344613414188
$ws.Range("J8").Value = -0.5052931892951553
$ws.Range("E9").Value = -0.6956293248699703
$ws.Range("J9").Value = 0.6842283227823049
$ws.Range("B10").Value = -6.350371468051042
$ws.Range("C10").Value = 2.449694784569639
$ws.Range("D10").Value = -2.734589956512321
$ws.Range("E10").Value = -1.391201525421976
$ws.Range("F10").Value = -8.902724452679363
$ws.Range("G10").Value = -0.2050953498393789
$ws.Range("H10").Value = 0.5052931892951553
$ws.Range("I10").Value = -0.6842283227823049
